# TimeSheet.xlsx update:
#  - Turn the "28/11" / "29/11" text labels in column A into real dates
#    (formatted as d-mmm), matching the existing time-tracking rows.
#  - Append two new rows for the latest day of work (01/12), describing the
#    new enemy spawning / AI behaviour / inventory work.
#  - Move the active selection below the new data and tidy the last
#    description text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column A: give the existing entries a proper date value + format ---
$ws.Range("A2:A5").NumberFormat = "d-mmm"
$ws.Range("A2").Value = "11/28/2017"
$ws.Range("A3").Value = "11/28/2017"
$ws.Range("A4").Value = "11/28/2017"
$ws.Range("A5").Value = "11/29/2017"

# --- Correct the description for the 29/11 entry ---
$ws.Range("C5").Value = "Restructuring character + basic sword implementation + camera movement"

# --- New rows for 01/12: enemy spawning, enemy behaviour, basic inventory ---
$ws.Range("A6:A7").NumberFormat = "d-mmm"
$ws.Range("A6").Value = "12/1/2017"
$ws.Range("A7").Value = "12/1/2017"

$ws.Range("B6:B7").NumberFormat = "h:mm"
$ws.Range("B6").Value = 0.03125
$ws.Range("B7").Value = 0.03125

$ws.Range("C6").Value = "Spawn enemies through manager + killable enemies"
$ws.Range("C7").Value = "Basic inventory + enemy behaviour"

# --- Move selection to the first empty row below the new data ---
$ws.Range("A8").Select()
